$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "data last updated" timestamp in A1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 27 de Marzo de 2020 a las 17:44"

# Rewrite the country data table (rows 4-205), which has been re-sorted
# by total cases (column B, descending) after the data refresh, and many
# country figures were updated.
$ws.Cells.Item(4,1).Value = "Estados Unidos"
$ws.Cells.Item(4,2).Value = 93300
$ws.Cells.Item(4,3).Value = 7865
$ws.Cells.Item(4,4).Value = 2424
$ws.Cells.Item(4,5).Value = 89492
$ws.Cells.Item(4,6).Value = 2432
$ws.Cells.Item(4,7).Value = 89
$ws.Cells.Item(4,8).Value = 1384
$ws.Cells.Item(5,1).Value = "China"
$ws.Cells.Item(5,2).Value = 81340
$ws.Cells.Item(5,3).Value = 55
$ws.Cells.Item(5,4).Value = 74588
$ws.Cells.Item(5,5).Value = 3460
$ws.Cells.Item(5,6).Value = 1034
$ws.Cells.Item(5,7).Value = 5
$ws.Cells.Item(5,8).Value = 3292
$ws.Cells.Item(6,1).Value = "Italia"
$ws.Cells.Item(6,2).Value = 80589
$ws.Cells.Item(6,3).Value = 0
$ws.Cells.Item(6,4).Value = 10361
$ws.Cells.Item(6,5).Value = 62013
$ws.Cells.Item(6,6).Value = 3612
$ws.Cells.Item(6,7).Value = 0
$ws.Cells.Item(6,8).Value = 8215
$ws.Cells.Item(7,1).Value = "España"
$ws.Cells.Item(7,2).Value = 64059
$ws.Cells.Item(7,3).Value = 6273
$ws.Cells.Item(7,4).Value = 9357
$ws.Cells.Item(7,5).Value = 49768
$ws.Cells.Item(7,6).Value = 4165
$ws.Cells.Item(7,7).Value = 569
$ws.Cells.Item(7,8).Value = 4934
$ws.Cells.Item(8,1).Value = "Alemania"
$ws.Cells.Item(8,2).Value = 49344
$ws.Cells.Item(8,3).Value = 5406
$ws.Cells.Item(8,4).Value = 5673
$ws.Cells.Item(8,5).Value = 43367
$ws.Cells.Item(8,6).Value = 23
$ws.Cells.Item(8,7).Value = 37
$ws.Cells.Item(8,8).Value = 304
$ws.Cells.Item(9,1).Value = "Iran"
$ws.Cells.Item(9,2).Value = 32332
$ws.Cells.Item(9,3).Value = 2926
$ws.Cells.Item(9,4).Value = 11133
$ws.Cells.Item(9,5).Value = 18821
$ws.Cells.Item(9,6).Value = 2893
$ws.Cells.Item(9,7).Value = 144
$ws.Cells.Item(9,8).Value = 2378
$ws.Cells.Item(10,1).Value = "Francia"
$ws.Cells.Item(10,2).Value = 29155
$ws.Cells.Item(10,3).Value = 0
$ws.Cells.Item(10,4).Value = 4948
$ws.Cells.Item(10,5).Value = 22511
$ws.Cells.Item(10,6).Value = 3375
$ws.Cells.Item(10,7).Value = 0
$ws.Cells.Item(10,8).Value = 1696
$ws.Cells.Item(11,1).Value = "Reino Unido"
$ws.Cells.Item(11,2).Value = 14543
$ws.Cells.Item(11,3).Value = 2885
$ws.Cells.Item(11,4).Value = 135
$ws.Cells.Item(11,5).Value = 13649
$ws.Cells.Item(11,6).Value = 163
$ws.Cells.Item(11,7).Value = 181
$ws.Cells.Item(11,8).Value = 759
$ws.Cells.Item(12,1).Value = "Suiza"
$ws.Cells.Item(12,2).Value = 12311
$ws.Cells.Item(12,3).Value = 500
$ws.Cells.Item(12,4).Value = 897
$ws.Cells.Item(12,5).Value = 11207
$ws.Cells.Item(12,6).Value = 203
$ws.Cells.Item(12,7).Value = 15
$ws.Cells.Item(12,8).Value = 207
$ws.Cells.Item(13,1).Value = "Corea del Sur"
$ws.Cells.Item(13,2).Value = 9332
$ws.Cells.Item(13,3).Value = 91
$ws.Cells.Item(13,4).Value = 4528
$ws.Cells.Item(13,5).Value = 4665
$ws.Cells.Item(13,6).Value = 59
$ws.Cells.Item(13,7).Value = 8
$ws.Cells.Item(13,8).Value = 139
$ws.Cells.Item(14,1).Value = "Paises Bajos"
$ws.Cells.Item(14,2).Value = 8603
$ws.Cells.Item(14,3).Value = 1172
$ws.Cells.Item(14,4).Value = 3
$ws.Cells.Item(14,5).Value = 8054
$ws.Cells.Item(14,6).Value = 761
$ws.Cells.Item(14,7).Value = 112
$ws.Cells.Item(14,8).Value = 546
$ws.Cells.Item(15,1).Value = "Austria"
$ws.Cells.Item(15,2).Value = 7557
$ws.Cells.Item(15,3).Value = 648
$ws.Cells.Item(15,4).Value = 225
$ws.Cells.Item(15,5).Value = 7274
$ws.Cells.Item(15,6).Value = 128
$ws.Cells.Item(15,7).Value = 9
$ws.Cells.Item(15,8).Value = 58
$ws.Cells.Item(16,1).Value = "Belgica"
$ws.Cells.Item(16,2).Value = 7284
$ws.Cells.Item(16,3).Value = 1049
$ws.Cells.Item(16,4).Value = 858
$ws.Cells.Item(16,5).Value = 6137
$ws.Cells.Item(16,6).Value = 690
$ws.Cells.Item(16,7).Value = 69
$ws.Cells.Item(16,8).Value = 289
$ws.Cells.Item(17,1).Value = "Portugal"
$ws.Cells.Item(17,2).Value = 4268
$ws.Cells.Item(17,3).Value = 724
$ws.Cells.Item(17,4).Value = 43
$ws.Cells.Item(17,5).Value = 4149
$ws.Cells.Item(17,6).Value = 71
$ws.Cells.Item(17,7).Value = 16
$ws.Cells.Item(17,8).Value = 76
$ws.Cells.Item(18,1).Value = "Canada"
$ws.Cells.Item(18,2).Value = 4043
$ws.Cells.Item(18,3).Value = 0
$ws.Cells.Item(18,4).Value = 228
$ws.Cells.Item(18,5).Value = 3776
$ws.Cells.Item(18,6).Value = 120
$ws.Cells.Item(18,7).Value = 0
$ws.Cells.Item(18,8).Value = 39
$ws.Cells.Item(19,1).Value = "Noruega"
$ws.Cells.Item(19,2).Value = 3696
$ws.Cells.Item(19,3).Value = 324
$ws.Cells.Item(19,4).Value = 6
$ws.Cells.Item(19,5).Value = 3671
$ws.Cells.Item(19,6).Value = 70
$ws.Cells.Item(19,7).Value = 5
$ws.Cells.Item(19,8).Value = 19
$ws.Cells.Item(20,1).Value = "Turquia"
$ws.Cells.Item(20,2).Value = 3629
$ws.Cells.Item(20,3).Value = 0
$ws.Cells.Item(20,4).Value = 26
$ws.Cells.Item(20,5).Value = 3528
$ws.Cells.Item(20,6).Value = 136
$ws.Cells.Item(20,7).Value = 0
$ws.Cells.Item(20,8).Value = 75
$ws.Cells.Item(21,1).Value = "Australia"
$ws.Cells.Item(21,2).Value = 3180
$ws.Cells.Item(21,3).Value = 130
$ws.Cells.Item(21,4).Value = 170
$ws.Cells.Item(21,5).Value = 2997
$ws.Cells.Item(21,6).Value = 23
$ws.Cells.Item(21,7).Value = 0
$ws.Cells.Item(21,8).Value = 13
$ws.Cells.Item(22,1).Value = "Suecia"
$ws.Cells.Item(22,2).Value = 3046
$ws.Cells.Item(22,3).Value = 206
$ws.Cells.Item(22,4).Value = 16
$ws.Cells.Item(22,5).Value = 2938
$ws.Cells.Item(22,6).Value = 214
$ws.Cells.Item(22,7).Value = 15
$ws.Cells.Item(22,8).Value = 92
$ws.Cells.Item(23,1).Value = "Israel"
$ws.Cells.Item(23,2).Value = 3035
$ws.Cells.Item(23,3).Value = 342
$ws.Cells.Item(23,4).Value = 79
$ws.Cells.Item(23,5).Value = 2944
$ws.Cells.Item(23,6).Value = 49
$ws.Cells.Item(23,7).Value = 4
$ws.Cells.Item(23,8).Value = 12
$ws.Cells.Item(24,1).Value = "Brasil"
$ws.Cells.Item(24,2).Value = 3027
$ws.Cells.Item(24,3).Value = 42
$ws.Cells.Item(24,4).Value = 6
$ws.Cells.Item(24,5).Value = 2944
$ws.Cells.Item(24,6).Value = 296
$ws.Cells.Item(24,7).Value = 0
$ws.Cells.Item(24,8).Value = 77
$ws.Cells.Item(25,1).Value = "Malasia"
$ws.Cells.Item(25,2).Value = 2161
$ws.Cells.Item(25,3).Value = 130
$ws.Cells.Item(25,4).Value = 259
$ws.Cells.Item(25,5).Value = 1876
$ws.Cells.Item(25,6).Value = 54
$ws.Cells.Item(25,7).Value = 3
$ws.Cells.Item(25,8).Value = 26
$ws.Cells.Item(26,1).Value = "Chequia"
$ws.Cells.Item(26,2).Value = 2062
$ws.Cells.Item(26,3).Value = 137
$ws.Cells.Item(26,4).Value = 11
$ws.Cells.Item(26,5).Value = 2042
$ws.Cells.Item(26,6).Value = 34
$ws.Cells.Item(26,7).Value = 0
$ws.Cells.Item(26,8).Value = 9
$ws.Cells.Item(27,1).Value = "Dinamarca"
$ws.Cells.Item(27,2).Value = 2010
$ws.Cells.Item(27,3).Value = 133
$ws.Cells.Item(27,4).Value = 1
$ws.Cells.Item(27,5).Value = 1957
$ws.Cells.Item(27,6).Value = 109
$ws.Cells.Item(27,7).Value = 11
$ws.Cells.Item(27,8).Value = 52
$ws.Cells.Item(28,1).Value = "Irlanda"
$ws.Cells.Item(28,2).Value = 1819
$ws.Cells.Item(28,3).Value = 0
$ws.Cells.Item(28,4).Value = 5
$ws.Cells.Item(28,5).Value = 1795
$ws.Cells.Item(28,6).Value = 47
$ws.Cells.Item(28,7).Value = 0
$ws.Cells.Item(28,8).Value = 19
$ws.Cells.Item(29,1).Value = "Chile"
$ws.Cells.Item(29,2).Value = 1610
$ws.Cells.Item(29,3).Value = 304
$ws.Cells.Item(29,4).Value = 43
$ws.Cells.Item(29,5).Value = 1562
$ws.Cells.Item(29,6).Value = 7
$ws.Cells.Item(29,7).Value = 1
$ws.Cells.Item(29,8).Value = 5
$ws.Cells.Item(30,1).Value = "Luxemburgo"
$ws.Cells.Item(30,2).Value = 1605
$ws.Cells.Item(30,3).Value = 152
$ws.Cells.Item(30,4).Value = 40
$ws.Cells.Item(30,5).Value = 1550
$ws.Cells.Item(30,6).Value = 25
$ws.Cells.Item(30,7).Value = 6
$ws.Cells.Item(30,8).Value = 15
$ws.Cells.Item(31,1).Value = "Ecuador"
$ws.Cells.Item(31,2).Value = 1595
$ws.Cells.Item(31,3).Value = 192
$ws.Cells.Item(31,4).Value = 3
$ws.Cells.Item(31,5).Value = 1556
$ws.Cells.Item(31,6).Value = 58
$ws.Cells.Item(31,7).Value = 2
$ws.Cells.Item(31,8).Value = 36
$ws.Cells.Item(32,1).Value = "Japon"
$ws.Cells.Item(32,2).Value = 1387
$ws.Cells.Item(32,3).Value = 0
$ws.Cells.Item(32,4).Value = 359
$ws.Cells.Item(32,5).Value = 981
$ws.Cells.Item(32,6).Value = 57
$ws.Cells.Item(32,7).Value = 0
$ws.Cells.Item(32,8).Value = 47
$ws.Cells.Item(33,1).Value = "Pakistan"
$ws.Cells.Item(33,2).Value = 1296
$ws.Cells.Item(33,3).Value = 95
$ws.Cells.Item(33,4).Value = 23
$ws.Cells.Item(33,5).Value = 1264
$ws.Cells.Item(33,6).Value = 7
$ws.Cells.Item(33,7).Value = 0
$ws.Cells.Item(33,8).Value = 9
$ws.Cells.Item(34,1).Value = "Rumania"
$ws.Cells.Item(34,2).Value = 1292
$ws.Cells.Item(34,3).Value = 263
$ws.Cells.Item(34,4).Value = 115
$ws.Cells.Item(34,5).Value = 1153
$ws.Cells.Item(34,6).Value = 32
$ws.Cells.Item(34,7).Value = 1
$ws.Cells.Item(34,8).Value = 24
$ws.Cells.Item(35,1).Value = "Polonia"
$ws.Cells.Item(35,2).Value = 1289
$ws.Cells.Item(35,3).Value = 68
$ws.Cells.Item(35,4).Value = 7
$ws.Cells.Item(35,5).Value = 1266
$ws.Cells.Item(35,6).Value = 3
$ws.Cells.Item(35,7).Value = 0
$ws.Cells.Item(35,8).Value = 16
$ws.Cells.Item(36,1).Value = "Tailandia"
$ws.Cells.Item(36,2).Value = 1136
$ws.Cells.Item(36,3).Value = 91
$ws.Cells.Item(36,4).Value = 97
$ws.Cells.Item(36,5).Value = 1034
$ws.Cells.Item(36,6).Value = 11
$ws.Cells.Item(36,7).Value = 1
$ws.Cells.Item(36,8).Value = 5
$ws.Cells.Item(37,1).Value = "Arabia Saudita"
$ws.Cells.Item(37,2).Value = 1104
$ws.Cells.Item(37,3).Value = 92
$ws.Cells.Item(37,4).Value = 35
$ws.Cells.Item(37,5).Value = 1066
$ws.Cells.Item(37,6).Value = 6
$ws.Cells.Item(37,7).Value = 0
$ws.Cells.Item(37,8).Value = 3
$ws.Cells.Item(38,1).Value = "Indonesia"
$ws.Cells.Item(38,2).Value = 1046
$ws.Cells.Item(38,3).Value = 153
$ws.Cells.Item(38,4).Value = 46
$ws.Cells.Item(38,5).Value = 913
$ws.Cells.Item(38,6).Value = 0
$ws.Cells.Item(38,7).Value = 9
$ws.Cells.Item(38,8).Value = 87
$ws.Cells.Item(39,1).Value = "Finlandia"
$ws.Cells.Item(39,2).Value = 1041
$ws.Cells.Item(39,3).Value = 83
$ws.Cells.Item(39,4).Value = 10
$ws.Cells.Item(39,5).Value = 1024
$ws.Cells.Item(39,6).Value = 32
$ws.Cells.Item(39,7).Value = 2
$ws.Cells.Item(39,8).Value = 7
$ws.Cells.Item(40,1).Value = "Rusia"
$ws.Cells.Item(40,2).Value = 1036
$ws.Cells.Item(40,3).Value = 196
$ws.Cells.Item(40,4).Value = 45
$ws.Cells.Item(40,5).Value = 988
$ws.Cells.Item(40,6).Value = 8
$ws.Cells.Item(40,7).Value = 0
$ws.Cells.Item(40,8).Value = 3
$ws.Cells.Item(41,1).Value = "Grecia"
$ws.Cells.Item(41,2).Value = 966
$ws.Cells.Item(41,3).Value = 74
$ws.Cells.Item(41,4).Value = 52
$ws.Cells.Item(41,5).Value = 886
$ws.Cells.Item(41,6).Value = 66
$ws.Cells.Item(41,7).Value = 1
$ws.Cells.Item(41,8).Value = 28
$ws.Cells.Item(42,1).Value = "Sudafrica"
$ws.Cells.Item(42,2).Value = 927
$ws.Cells.Item(42,3).Value = 0
$ws.Cells.Item(42,4).Value = 12
$ws.Cells.Item(42,5).Value = 913
$ws.Cells.Item(42,6).Value = 7
$ws.Cells.Item(42,7).Value = 2
$ws.Cells.Item(42,8).Value = 2
$ws.Cells.Item(43,1).Value = "Islandia"
$ws.Cells.Item(43,2).Value = 890
$ws.Cells.Item(43,3).Value = 88
$ws.Cells.Item(43,4).Value = 97
$ws.Cells.Item(43,5).Value = 791
$ws.Cells.Item(43,6).Value = 18
$ws.Cells.Item(43,7).Value = 0
$ws.Cells.Item(43,8).Value = 2
$ws.Cells.Item(44,1).Value = "India"
$ws.Cells.Item(44,2).Value = 878
$ws.Cells.Item(44,3).Value = 151
$ws.Cells.Item(44,4).Value = 73
$ws.Cells.Item(44,5).Value = 785
$ws.Cells.Item(44,6).Value = 0
$ws.Cells.Item(44,7).Value = 0
$ws.Cells.Item(44,8).Value = 20
$ws.Cells.Item(45,1).Value = "Filipinas"
$ws.Cells.Item(45,2).Value = 803
$ws.Cells.Item(45,3).Value = 96
$ws.Cells.Item(45,4).Value = 31
$ws.Cells.Item(45,5).Value = 718
$ws.Cells.Item(45,6).Value = 1
$ws.Cells.Item(45,7).Value = 9
$ws.Cells.Item(45,8).Value = 54
$ws.Cells.Item(46,1).Value = "Singapur"
$ws.Cells.Item(46,2).Value = 732
$ws.Cells.Item(46,3).Value = 49
$ws.Cells.Item(46,4).Value = 183
$ws.Cells.Item(46,5).Value = 547
$ws.Cells.Item(46,6).Value = 17
$ws.Cells.Item(46,7).Value = 0
$ws.Cells.Item(46,8).Value = 2
$ws.Cells.Item(47,1).Value = "Crucero"
$ws.Cells.Item(47,2).Value = 712
$ws.Cells.Item(47,3).Value = 0
$ws.Cells.Item(47,4).Value = 597
$ws.Cells.Item(47,5).Value = 105
$ws.Cells.Item(47,6).Value = 15
$ws.Cells.Item(47,7).Value = 0
$ws.Cells.Item(47,8).Value = 10
$ws.Cells.Item(48,1).Value = "Panama"
$ws.Cells.Item(48,2).Value = 674
$ws.Cells.Item(48,3).Value = 0
$ws.Cells.Item(48,4).Value = 2
$ws.Cells.Item(48,5).Value = 663
$ws.Cells.Item(48,6).Value = 20
$ws.Cells.Item(48,7).Value = 0
$ws.Cells.Item(48,8).Value = 9
$ws.Cells.Item(49,1).Value = "Eslovenia"
$ws.Cells.Item(49,2).Value = 632
$ws.Cells.Item(49,3).Value = 70
$ws.Cells.Item(49,4).Value = 10
$ws.Cells.Item(49,5).Value = 613
$ws.Cells.Item(49,6).Value = 14
$ws.Cells.Item(49,7).Value = 3
$ws.Cells.Item(49,8).Value = 9
$ws.Cells.Item(50,1).Value = "Argentina"
$ws.Cells.Item(50,2).Value = 589
$ws.Cells.Item(50,3).Value = 0
$ws.Cells.Item(50,4).Value = 72
$ws.Cells.Item(50,5).Value = 504
$ws.Cells.Item(50,6).Value = 0
$ws.Cells.Item(50,7).Value = 1
$ws.Cells.Item(50,8).Value = 13
$ws.Cells.Item(51,1).Value = "Croacia"
$ws.Cells.Item(51,2).Value = 586
$ws.Cells.Item(51,3).Value = 91
$ws.Cells.Item(51,4).Value = 37
$ws.Cells.Item(51,5).Value = 546
$ws.Cells.Item(51,6).Value = 14
$ws.Cells.Item(51,7).Value = 0
$ws.Cells.Item(51,8).Value = 3
$ws.Cells.Item(52,1).Value = "Mexico"
$ws.Cells.Item(52,2).Value = 585
$ws.Cells.Item(52,3).Value = 110
$ws.Cells.Item(52,4).Value = 4
$ws.Cells.Item(52,5).Value = 573
$ws.Cells.Item(52,6).Value = 1
$ws.Cells.Item(52,7).Value = 2
$ws.Cells.Item(52,8).Value = 8
$ws.Cells.Item(53,1).Value = "Republica Dominicana"
$ws.Cells.Item(53,2).Value = 581
$ws.Cells.Item(53,3).Value = 93
$ws.Cells.Item(53,4).Value = 3
$ws.Cells.Item(53,5).Value = 558
$ws.Cells.Item(53,6).Value = 0
$ws.Cells.Item(53,7).Value = 10
$ws.Cells.Item(53,8).Value = 20
$ws.Cells.Item(54,1).Value = "Peru"
$ws.Cells.Item(54,2).Value = 580
$ws.Cells.Item(54,3).Value = 0
$ws.Cells.Item(54,4).Value = 14
$ws.Cells.Item(54,5).Value = 557
$ws.Cells.Item(54,6).Value = 14
$ws.Cells.Item(54,7).Value = 0
$ws.Cells.Item(54,8).Value = 9
$ws.Cells.Item(55,1).Value = "Estonia"
$ws.Cells.Item(55,2).Value = 575
$ws.Cells.Item(55,3).Value = 37
$ws.Cells.Item(55,4).Value = 11
$ws.Cells.Item(55,5).Value = 563
$ws.Cells.Item(55,6).Value = 7
$ws.Cells.Item(55,7).Value = 0
$ws.Cells.Item(55,8).Value = 1
$ws.Cells.Item(56,1).Value = "Catar"
$ws.Cells.Item(56,2).Value = 562
$ws.Cells.Item(56,3).Value = 13
$ws.Cells.Item(56,4).Value = 43
$ws.Cells.Item(56,5).Value = 519
$ws.Cells.Item(56,6).Value = 6
$ws.Cells.Item(56,7).Value = 0
$ws.Cells.Item(56,8).Value = 0
$ws.Cells.Item(57,1).Value = "Serbia"
$ws.Cells.Item(57,2).Value = 528
$ws.Cells.Item(57,3).Value = 71
$ws.Cells.Item(57,4).Value = 15
$ws.Cells.Item(57,5).Value = 505
$ws.Cells.Item(57,6).Value = 25
$ws.Cells.Item(57,7).Value = 1
$ws.Cells.Item(57,8).Value = 8
$ws.Cells.Item(58,1).Value = "Hong Kong"
$ws.Cells.Item(58,2).Value = 518
$ws.Cells.Item(58,3).Value = 64
$ws.Cells.Item(58,4).Value = 111
$ws.Cells.Item(58,5).Value = 403
$ws.Cells.Item(58,6).Value = 5
$ws.Cells.Item(58,7).Value = 0
$ws.Cells.Item(58,8).Value = 4
$ws.Cells.Item(59,1).Value = "Egipto"
$ws.Cells.Item(59,2).Value = 495
$ws.Cells.Item(59,3).Value = 0
$ws.Cells.Item(59,4).Value = 102
$ws.Cells.Item(59,5).Value = 369
$ws.Cells.Item(59,6).Value = 0
$ws.Cells.Item(59,7).Value = 0
$ws.Cells.Item(59,8).Value = 24
$ws.Cells.Item(60,1).Value = "Colombia"
$ws.Cells.Item(60,2).Value = 491
$ws.Cells.Item(60,3).Value = 0
$ws.Cells.Item(60,4).Value = 8
$ws.Cells.Item(60,5).Value = 477
$ws.Cells.Item(60,6).Value = 0
$ws.Cells.Item(60,7).Value = 0
$ws.Cells.Item(60,8).Value = 6
$ws.Cells.Item(61,1).Value = "Barein"
$ws.Cells.Item(61,2).Value = 466
$ws.Cells.Item(61,3).Value = 8
$ws.Cells.Item(61,4).Value = 227
$ws.Cells.Item(61,5).Value = 235
$ws.Cells.Item(61,6).Value = 1
$ws.Cells.Item(61,7).Value = 0
$ws.Cells.Item(61,8).Value = 4
$ws.Cells.Item(62,1).Value = "Irak"
$ws.Cells.Item(62,2).Value = 458
$ws.Cells.Item(62,3).Value = 76
$ws.Cells.Item(62,4).Value = 122
$ws.Cells.Item(62,5).Value = 296
$ws.Cells.Item(62,6).Value = 0
$ws.Cells.Item(62,7).Value = 4
$ws.Cells.Item(62,8).Value = 40
$ws.Cells.Item(63,1).Value = "Argelia"
$ws.Cells.Item(63,2).Value = 409
$ws.Cells.Item(63,3).Value = 42
$ws.Cells.Item(63,4).Value = 29
$ws.Cells.Item(63,5).Value = 354
$ws.Cells.Item(63,6).Value = 0
$ws.Cells.Item(63,7).Value = 1
$ws.Cells.Item(63,8).Value = 26
$ws.Cells.Item(64,1).Value = "Libano"
$ws.Cells.Item(64,2).Value = 391
$ws.Cells.Item(64,3).Value = 23
$ws.Cells.Item(64,4).Value = 23
$ws.Cells.Item(64,5).Value = 361
$ws.Cells.Item(64,6).Value = 3
$ws.Cells.Item(64,7).Value = 1
$ws.Cells.Item(64,8).Value = 7
$ws.Cells.Item(65,1).Value = "Nueva Zelanda"
$ws.Cells.Item(65,2).Value = 368
$ws.Cells.Item(65,3).Value = 0
$ws.Cells.Item(65,4).Value = 37
$ws.Cells.Item(65,5).Value = 331
$ws.Cells.Item(65,6).Value = 1
$ws.Cells.Item(65,7).Value = 0
$ws.Cells.Item(65,8).Value = 0
$ws.Cells.Item(66,1).Value = "Lituania"
$ws.Cells.Item(66,2).Value = 345
$ws.Cells.Item(66,3).Value = 46
$ws.Cells.Item(66,4).Value = 1
$ws.Cells.Item(66,5).Value = 339
$ws.Cells.Item(66,6).Value = 2
$ws.Cells.Item(66,7).Value = 1
$ws.Cells.Item(66,8).Value = 5
$ws.Cells.Item(67,1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(67,2).Value = 333
$ws.Cells.Item(67,3).Value = 0
$ws.Cells.Item(67,4).Value = 52
$ws.Cells.Item(67,5).Value = 279
$ws.Cells.Item(67,6).Value = 2
$ws.Cells.Item(67,7).Value = 0
$ws.Cells.Item(67,8).Value = 2
$ws.Cells.Item(68,1).Value = "Armenia"
$ws.Cells.Item(68,2).Value = 329
$ws.Cells.Item(68,3).Value = 39
$ws.Cells.Item(68,4).Value = 28
$ws.Cells.Item(68,5).Value = 300
$ws.Cells.Item(68,6).Value = 6
$ws.Cells.Item(68,7).Value = 0
$ws.Cells.Item(68,8).Value = 1
$ws.Cells.Item(69,1).Value = "Hungria"
$ws.Cells.Item(69,2).Value = 300
$ws.Cells.Item(69,3).Value = 39
$ws.Cells.Item(69,4).Value = 34
$ws.Cells.Item(69,5).Value = 256
$ws.Cells.Item(69,6).Value = 6
$ws.Cells.Item(69,7).Value = 0
$ws.Cells.Item(69,8).Value = 10
$ws.Cells.Item(70,1).Value = "Bulgaria"
$ws.Cells.Item(70,2).Value = 293
$ws.Cells.Item(70,3).Value = 29
$ws.Cells.Item(70,4).Value = 9
$ws.Cells.Item(70,5).Value = 281
$ws.Cells.Item(70,6).Value = 8
$ws.Cells.Item(70,7).Value = 0
$ws.Cells.Item(70,8).Value = 3
$ws.Cells.Item(71,1).Value = "Letonia"
$ws.Cells.Item(71,2).Value = 280
$ws.Cells.Item(71,3).Value = 36
$ws.Cells.Item(71,4).Value = 1
$ws.Cells.Item(71,5).Value = 279
$ws.Cells.Item(71,6).Value = 0
$ws.Cells.Item(71,7).Value = 0
$ws.Cells.Item(71,8).Value = 0
$ws.Cells.Item(72,1).Value = "Marruecos"
$ws.Cells.Item(72,2).Value = 275
$ws.Cells.Item(72,3).Value = 0
$ws.Cells.Item(72,4).Value = 8
$ws.Cells.Item(72,5).Value = 256
$ws.Cells.Item(72,6).Value = 1
$ws.Cells.Item(72,7).Value = 0
$ws.Cells.Item(72,8).Value = 11
$ws.Cells.Item(73,1).Value = "Eslovaquia"
$ws.Cells.Item(73,2).Value = 269
$ws.Cells.Item(73,3).Value = 43
$ws.Cells.Item(73,4).Value = 2
$ws.Cells.Item(73,5).Value = 267
$ws.Cells.Item(73,6).Value = 1
$ws.Cells.Item(73,7).Value = 0
$ws.Cells.Item(73,8).Value = 0
$ws.Cells.Item(74,1).Value = "Principado de Andorra"
$ws.Cells.Item(74,2).Value = 267
$ws.Cells.Item(74,3).Value = 43
$ws.Cells.Item(74,4).Value = 1
$ws.Cells.Item(74,5).Value = 263
$ws.Cells.Item(74,6).Value = 11
$ws.Cells.Item(74,7).Value = 0
$ws.Cells.Item(74,8).Value = 3
$ws.Cells.Item(75,1).Value = "Taiwan"
$ws.Cells.Item(75,2).Value = 267
$ws.Cells.Item(75,3).Value = 15
$ws.Cells.Item(75,4).Value = 30
$ws.Cells.Item(75,5).Value = 235
$ws.Cells.Item(75,6).Value = 0
$ws.Cells.Item(75,7).Value = 0
$ws.Cells.Item(75,8).Value = 2
$ws.Cells.Item(76,1).Value = "Uruguay"
$ws.Cells.Item(76,2).Value = 238
$ws.Cells.Item(76,3).Value = 0
$ws.Cells.Item(76,4).Value = 0
$ws.Cells.Item(76,5).Value = 238
$ws.Cells.Item(76,6).Value = 3
$ws.Cells.Item(76,7).Value = 0
$ws.Cells.Item(76,8).Value = 0
$ws.Cells.Item(77,1).Value = "Costa Rica"
$ws.Cells.Item(77,2).Value = 231
$ws.Cells.Item(77,3).Value = 0
$ws.Cells.Item(77,4).Value = 2
$ws.Cells.Item(77,5).Value = 227
$ws.Cells.Item(77,6).Value = 5
$ws.Cells.Item(77,7).Value = 0
$ws.Cells.Item(77,8).Value = 2
$ws.Cells.Item(78,1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(78,2).Value = 231
$ws.Cells.Item(78,3).Value = 40
$ws.Cells.Item(78,4).Value = 5
$ws.Cells.Item(78,5).Value = 222
$ws.Cells.Item(78,6).Value = 1
$ws.Cells.Item(78,7).Value = 1
$ws.Cells.Item(78,8).Value = 4
$ws.Cells.Item(79,1).Value = "Tunez"
$ws.Cells.Item(79,2).Value = 227
$ws.Cells.Item(79,3).Value = 30
$ws.Cells.Item(79,4).Value = 2
$ws.Cells.Item(79,5).Value = 219
$ws.Cells.Item(79,6).Value = 10
$ws.Cells.Item(79,7).Value = 1
$ws.Cells.Item(79,8).Value = 6
$ws.Cells.Item(80,1).Value = "Ucrania"
$ws.Cells.Item(80,2).Value = 226
$ws.Cells.Item(80,3).Value = 30
$ws.Cells.Item(80,4).Value = 5
$ws.Cells.Item(80,5).Value = 216
$ws.Cells.Item(80,6).Value = 0
$ws.Cells.Item(80,7).Value = 0
$ws.Cells.Item(80,8).Value = 5
$ws.Cells.Item(81,1).Value = "Kuwait"
$ws.Cells.Item(81,2).Value = 225
$ws.Cells.Item(81,3).Value = 17
$ws.Cells.Item(81,4).Value = 57
$ws.Cells.Item(81,5).Value = 168
$ws.Cells.Item(81,6).Value = 11
$ws.Cells.Item(81,7).Value = 0
$ws.Cells.Item(81,8).Value = 0
$ws.Cells.Item(82,1).Value = "San Marino"
$ws.Cells.Item(82,2).Value = 223
$ws.Cells.Item(82,3).Value = 15
$ws.Cells.Item(82,4).Value = 4
$ws.Cells.Item(82,5).Value = 198
$ws.Cells.Item(82,6).Value = 12
$ws.Cells.Item(82,7).Value = 0
$ws.Cells.Item(82,8).Value = 21
$ws.Cells.Item(83,1).Value = "Republica de Macedonia"
$ws.Cells.Item(83,2).Value = 219
$ws.Cells.Item(83,3).Value = 18
$ws.Cells.Item(83,4).Value = 3
$ws.Cells.Item(83,5).Value = 213
$ws.Cells.Item(83,6).Value = 1
$ws.Cells.Item(83,7).Value = 0
$ws.Cells.Item(83,8).Value = 3
$ws.Cells.Item(84,1).Value = "Jordania"
$ws.Cells.Item(84,2).Value = 212
$ws.Cells.Item(84,3).Value = 0
$ws.Cells.Item(84,4).Value = 2
$ws.Cells.Item(84,5).Value = 210
$ws.Cells.Item(84,6).Value = 0
$ws.Cells.Item(84,7).Value = 0
$ws.Cells.Item(84,8).Value = 0
$ws.Cells.Item(85,1).Value = "Moldavia"
$ws.Cells.Item(85,2).Value = 199
$ws.Cells.Item(85,3).Value = 22
$ws.Cells.Item(85,4).Value = 2
$ws.Cells.Item(85,5).Value = 195
$ws.Cells.Item(85,6).Value = 33
$ws.Cells.Item(85,7).Value = 1
$ws.Cells.Item(85,8).Value = 2
$ws.Cells.Item(86,1).Value = "Albania"
$ws.Cells.Item(86,2).Value = 186
$ws.Cells.Item(86,3).Value = 12
$ws.Cells.Item(86,4).Value = 31
$ws.Cells.Item(86,5).Value = 147
$ws.Cells.Item(86,6).Value = 3
$ws.Cells.Item(86,7).Value = 2
$ws.Cells.Item(86,8).Value = 8
$ws.Cells.Item(87,1).Value = "Burkina Faso"
$ws.Cells.Item(87,2).Value = 180
$ws.Cells.Item(87,3).Value = 28
$ws.Cells.Item(87,4).Value = 12
$ws.Cells.Item(87,5).Value = 159
$ws.Cells.Item(87,6).Value = 0
$ws.Cells.Item(87,7).Value = 2
$ws.Cells.Item(87,8).Value = 9
$ws.Cells.Item(88,1).Value = "Azerbaiyan"
$ws.Cells.Item(88,2).Value = 165
$ws.Cells.Item(88,3).Value = 43
$ws.Cells.Item(88,4).Value = 15
$ws.Cells.Item(88,5).Value = 147
$ws.Cells.Item(88,6).Value = 6
$ws.Cells.Item(88,7).Value = 0
$ws.Cells.Item(88,8).Value = 3
$ws.Cells.Item(89,1).Value = "Vietnam"
$ws.Cells.Item(89,2).Value = 163
$ws.Cells.Item(89,3).Value = 10
$ws.Cells.Item(89,4).Value = 20
$ws.Cells.Item(89,5).Value = 143
$ws.Cells.Item(89,6).Value = 3
$ws.Cells.Item(89,7).Value = 0
$ws.Cells.Item(89,8).Value = 0
$ws.Cells.Item(90,1).Value = "Republica de Chipre"
$ws.Cells.Item(90,2).Value = 146
$ws.Cells.Item(90,3).Value = 0
$ws.Cells.Item(90,4).Value = 4
$ws.Cells.Item(90,5).Value = 139
$ws.Cells.Item(90,6).Value = 3
$ws.Cells.Item(90,7).Value = 0
$ws.Cells.Item(90,8).Value = 3
$ws.Cells.Item(91,1).Value = "Reunion"
$ws.Cells.Item(91,2).Value = 145
$ws.Cells.Item(91,3).Value = 10
$ws.Cells.Item(91,4).Value = 1
$ws.Cells.Item(91,5).Value = 144
$ws.Cells.Item(91,6).Value = 0
$ws.Cells.Item(91,7).Value = 0
$ws.Cells.Item(91,8).Value = 0
$ws.Cells.Item(92,1).Value = "Islas Feroe"
$ws.Cells.Item(92,2).Value = 144
$ws.Cells.Item(92,3).Value = 4
$ws.Cells.Item(92,4).Value = 54
$ws.Cells.Item(92,5).Value = 90
$ws.Cells.Item(92,6).Value = 2
$ws.Cells.Item(92,7).Value = 0
$ws.Cells.Item(92,8).Value = 0
$ws.Cells.Item(93,1).Value = "Malta"
$ws.Cells.Item(93,2).Value = 139
$ws.Cells.Item(93,3).Value = 5
$ws.Cells.Item(93,4).Value = 2
$ws.Cells.Item(93,5).Value = 137
$ws.Cells.Item(93,6).Value = 1
$ws.Cells.Item(93,7).Value = 0
$ws.Cells.Item(93,8).Value = 0
$ws.Cells.Item(94,1).Value = "Kazajistan"
$ws.Cells.Item(94,2).Value = 137
$ws.Cells.Item(94,3).Value = 24
$ws.Cells.Item(94,4).Value = 3
$ws.Cells.Item(94,5).Value = 133
$ws.Cells.Item(94,6).Value = 0
$ws.Cells.Item(94,7).Value = 0
$ws.Cells.Item(94,8).Value = 1
$ws.Cells.Item(95,1).Value = "Ghana"
$ws.Cells.Item(95,2).Value = 136
$ws.Cells.Item(95,3).Value = 4
$ws.Cells.Item(95,4).Value = 1
$ws.Cells.Item(95,5).Value = 131
$ws.Cells.Item(95,6).Value = 1
$ws.Cells.Item(95,7).Value = 0
$ws.Cells.Item(95,8).Value = 4
$ws.Cells.Item(96,1).Value = "Oman"
$ws.Cells.Item(96,2).Value = 131
$ws.Cells.Item(96,3).Value = 22
$ws.Cells.Item(96,4).Value = 23
$ws.Cells.Item(96,5).Value = 108
$ws.Cells.Item(96,6).Value = 0
$ws.Cells.Item(96,7).Value = 0
$ws.Cells.Item(96,8).Value = 0
$ws.Cells.Item(97,1).Value = "Senegal"
$ws.Cells.Item(97,2).Value = 119
$ws.Cells.Item(97,3).Value = 14
$ws.Cells.Item(97,4).Value = 11
$ws.Cells.Item(97,5).Value = 108
$ws.Cells.Item(97,6).Value = 0
$ws.Cells.Item(97,7).Value = 0
$ws.Cells.Item(97,8).Value = 0
$ws.Cells.Item(98,1).Value = "Brunei"
$ws.Cells.Item(98,2).Value = 115
$ws.Cells.Item(98,3).Value = 1
$ws.Cells.Item(98,4).Value = 11
$ws.Cells.Item(98,5).Value = 104
$ws.Cells.Item(98,6).Value = 1
$ws.Cells.Item(98,7).Value = 0
$ws.Cells.Item(98,8).Value = 0
$ws.Cells.Item(99,1).Value = "Venezuela"
$ws.Cells.Item(99,2).Value = 107
$ws.Cells.Item(99,3).Value = 0
$ws.Cells.Item(99,4).Value = 31
$ws.Cells.Item(99,5).Value = 75
$ws.Cells.Item(99,6).Value = 2
$ws.Cells.Item(99,7).Value = 0
$ws.Cells.Item(99,8).Value = 1
$ws.Cells.Item(100,1).Value = "Sri Lanka"
$ws.Cells.Item(100,2).Value = 106
$ws.Cells.Item(100,3).Value = 0
$ws.Cells.Item(100,4).Value = 7
$ws.Cells.Item(100,5).Value = 99
$ws.Cells.Item(100,6).Value = 5
$ws.Cells.Item(100,7).Value = 0
$ws.Cells.Item(100,8).Value = 0
$ws.Cells.Item(101,1).Value = "Camboya"
$ws.Cells.Item(101,2).Value = 99
$ws.Cells.Item(101,3).Value = 1
$ws.Cells.Item(101,4).Value = 11
$ws.Cells.Item(101,5).Value = 88
$ws.Cells.Item(101,6).Value = 1
$ws.Cells.Item(101,7).Value = 0
$ws.Cells.Item(101,8).Value = 0
$ws.Cells.Item(102,1).Value = "Costa de Marfil"
$ws.Cells.Item(102,2).Value = 96
$ws.Cells.Item(102,3).Value = 0
$ws.Cells.Item(102,4).Value = 3
$ws.Cells.Item(102,5).Value = 93
$ws.Cells.Item(102,6).Value = 0
$ws.Cells.Item(102,7).Value = 0
$ws.Cells.Item(102,8).Value = 0
$ws.Cells.Item(103,1).Value = "Mauricio"
$ws.Cells.Item(103,2).Value = 94
$ws.Cells.Item(103,3).Value = 13
$ws.Cells.Item(103,4).Value = 0
$ws.Cells.Item(103,5).Value = 92
$ws.Cells.Item(103,6).Value = 1
$ws.Cells.Item(103,7).Value = 0
$ws.Cells.Item(103,8).Value = 2
$ws.Cells.Item(104,1).Value = "Afganistan"
$ws.Cells.Item(104,2).Value = 94
$ws.Cells.Item(104,3).Value = 0
$ws.Cells.Item(104,4).Value = 2
$ws.Cells.Item(104,5).Value = 88
$ws.Cells.Item(104,6).Value = 0
$ws.Cells.Item(104,7).Value = 0
$ws.Cells.Item(104,8).Value = 4
$ws.Cells.Item(105,1).Value = "Bielorrusia"
$ws.Cells.Item(105,2).Value = 94
$ws.Cells.Item(105,3).Value = 8
$ws.Cells.Item(105,4).Value = 32
$ws.Cells.Item(105,5).Value = 62
$ws.Cells.Item(105,6).Value = 2
$ws.Cells.Item(105,7).Value = 0
$ws.Cells.Item(105,8).Value = 0
$ws.Cells.Item(106,1).Value = "Estado de Palestina"
$ws.Cells.Item(106,2).Value = 91
$ws.Cells.Item(106,3).Value = 5
$ws.Cells.Item(106,4).Value = 17
$ws.Cells.Item(106,5).Value = 73
$ws.Cells.Item(106,6).Value = 0
$ws.Cells.Item(106,7).Value = 0
$ws.Cells.Item(106,8).Value = 1
$ws.Cells.Item(107,1).Value = "Camerun"
$ws.Cells.Item(107,2).Value = 88
$ws.Cells.Item(107,3).Value = 13
$ws.Cells.Item(107,4).Value = 2
$ws.Cells.Item(107,5).Value = 84
$ws.Cells.Item(107,6).Value = 0
$ws.Cells.Item(107,7).Value = 1
$ws.Cells.Item(107,8).Value = 2
$ws.Cells.Item(108,1).Value = "Uzbekistan"
$ws.Cells.Item(108,2).Value = 88
$ws.Cells.Item(108,3).Value = 13
$ws.Cells.Item(108,4).Value = 5
$ws.Cells.Item(108,5).Value = 82
$ws.Cells.Item(108,6).Value = 8
$ws.Cells.Item(108,7).Value = 1
$ws.Cells.Item(108,8).Value = 1
$ws.Cells.Item(109,1).Value = "Martinica"
$ws.Cells.Item(109,2).Value = 81
$ws.Cells.Item(109,3).Value = 0
$ws.Cells.Item(109,4).Value = 0
$ws.Cells.Item(109,5).Value = 80
$ws.Cells.Item(109,6).Value = 12
$ws.Cells.Item(109,7).Value = 0
$ws.Cells.Item(109,8).Value = 1
$ws.Cells.Item(110,1).Value = "Georgia"
$ws.Cells.Item(110,2).Value = 81
$ws.Cells.Item(110,3).Value = 2
$ws.Cells.Item(110,4).Value = 13
$ws.Cells.Item(110,5).Value = 68
$ws.Cells.Item(110,6).Value = 1
$ws.Cells.Item(110,7).Value = 0
$ws.Cells.Item(110,8).Value = 0
$ws.Cells.Item(111,1).Value = "Cuba"
$ws.Cells.Item(111,2).Value = 80
$ws.Cells.Item(111,3).Value = 13
$ws.Cells.Item(111,4).Value = 4
$ws.Cells.Item(111,5).Value = 74
$ws.Cells.Item(111,6).Value = 2
$ws.Cells.Item(111,7).Value = 0
$ws.Cells.Item(111,8).Value = 2
$ws.Cells.Item(112,1).Value = "Guadalupe"
$ws.Cells.Item(112,2).Value = 73
$ws.Cells.Item(112,3).Value = 0
$ws.Cells.Item(112,4).Value = 0
$ws.Cells.Item(112,5).Value = 72
$ws.Cells.Item(112,6).Value = 4
$ws.Cells.Item(112,7).Value = 0
$ws.Cells.Item(112,8).Value = 1
$ws.Cells.Item(113,1).Value = "Montenegro"
$ws.Cells.Item(113,2).Value = 70
$ws.Cells.Item(113,3).Value = 1
$ws.Cells.Item(113,4).Value = 0
$ws.Cells.Item(113,5).Value = 69
$ws.Cells.Item(113,6).Value = 1
$ws.Cells.Item(113,7).Value = 0
$ws.Cells.Item(113,8).Value = 1
$ws.Cells.Item(114,1).Value = "Honduras"
$ws.Cells.Item(114,2).Value = 68
$ws.Cells.Item(114,3).Value = 1
$ws.Cells.Item(114,4).Value = 0
$ws.Cells.Item(114,5).Value = 67
$ws.Cells.Item(114,6).Value = 0
$ws.Cells.Item(114,7).Value = 0
$ws.Cells.Item(114,8).Value = 1
$ws.Cells.Item(115,1).Value = "Trinidad yTobago"
$ws.Cells.Item(115,2).Value = 66
$ws.Cells.Item(115,3).Value = 1
$ws.Cells.Item(115,4).Value = 1
$ws.Cells.Item(115,5).Value = 63
$ws.Cells.Item(115,6).Value = 0
$ws.Cells.Item(115,7).Value = 1
$ws.Cells.Item(115,8).Value = 2
$ws.Cells.Item(116,1).Value = "Nigeria"
$ws.Cells.Item(116,2).Value = 65
$ws.Cells.Item(116,3).Value = 0
$ws.Cells.Item(116,4).Value = 3
$ws.Cells.Item(116,5).Value = 61
$ws.Cells.Item(116,6).Value = 0
$ws.Cells.Item(116,7).Value = 0
$ws.Cells.Item(116,8).Value = 1
$ws.Cells.Item(117,1).Value = "Bolivia"
$ws.Cells.Item(117,2).Value = 61
$ws.Cells.Item(117,3).Value = 0
$ws.Cells.Item(117,4).Value = 0
$ws.Cells.Item(117,5).Value = 61
$ws.Cells.Item(117,6).Value = 0
$ws.Cells.Item(117,7).Value = 0
$ws.Cells.Item(117,8).Value = 0
$ws.Cells.Item(118,1).Value = "Kirguistan"
$ws.Cells.Item(118,2).Value = 58
$ws.Cells.Item(118,3).Value = 14
$ws.Cells.Item(118,4).Value = 0
$ws.Cells.Item(118,5).Value = 58
$ws.Cells.Item(118,6).Value = 0
$ws.Cells.Item(118,7).Value = 0
$ws.Cells.Item(118,8).Value = 0
$ws.Cells.Item(119,1).Value = "Liechtenstein"
$ws.Cells.Item(119,2).Value = 56
$ws.Cells.Item(119,3).Value = 0
$ws.Cells.Item(119,4).Value = 0
$ws.Cells.Item(119,5).Value = 56
$ws.Cells.Item(119,6).Value = 0
$ws.Cells.Item(119,7).Value = 0
$ws.Cells.Item(119,8).Value = 0
$ws.Cells.Item(120,1).Value = "Gibraltar"
$ws.Cells.Item(120,2).Value = 55
$ws.Cells.Item(120,3).Value = 20
$ws.Cells.Item(120,4).Value = 14
$ws.Cells.Item(120,5).Value = 41
$ws.Cells.Item(120,6).Value = 0
$ws.Cells.Item(120,7).Value = 0
$ws.Cells.Item(120,8).Value = 0
$ws.Cells.Item(121,1).Value = "Paraguay"
$ws.Cells.Item(121,2).Value = 52
$ws.Cells.Item(121,3).Value = 11
$ws.Cells.Item(121,4).Value = 1
$ws.Cells.Item(121,5).Value = 48
$ws.Cells.Item(121,6).Value = 1
$ws.Cells.Item(121,7).Value = 0
$ws.Cells.Item(121,8).Value = 3
$ws.Cells.Item(122,1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(122,2).Value = 51
$ws.Cells.Item(122,3).Value = 0
$ws.Cells.Item(122,4).Value = 2
$ws.Cells.Item(122,5).Value = 46
$ws.Cells.Item(122,6).Value = 0
$ws.Cells.Item(122,7).Value = 0
$ws.Cells.Item(122,8).Value = 3
$ws.Cells.Item(123,1).Value = "Mayotte"
$ws.Cells.Item(123,2).Value = 50
$ws.Cells.Item(123,3).Value = 14
$ws.Cells.Item(123,4).Value = 0
$ws.Cells.Item(123,5).Value = 50
$ws.Cells.Item(123,6).Value = 0
$ws.Cells.Item(123,7).Value = 0
$ws.Cells.Item(123,8).Value = 0
$ws.Cells.Item(124,1).Value = "Ruanda"
$ws.Cells.Item(124,2).Value = 50
$ws.Cells.Item(124,3).Value = 0
$ws.Cells.Item(124,4).Value = 0
$ws.Cells.Item(124,5).Value = 50
$ws.Cells.Item(124,6).Value = 0
$ws.Cells.Item(124,7).Value = 0
$ws.Cells.Item(124,8).Value = 0
$ws.Cells.Item(125,1).Value = "Banglades"
$ws.Cells.Item(125,2).Value = 48
$ws.Cells.Item(125,3).Value = 4
$ws.Cells.Item(125,4).Value = 11
$ws.Cells.Item(125,5).Value = 32
$ws.Cells.Item(125,6).Value = 1
$ws.Cells.Item(125,7).Value = 0
$ws.Cells.Item(125,8).Value = 5
$ws.Cells.Item(126,1).Value = "Puerto Rico"
$ws.Cells.Item(126,2).Value = 39
$ws.Cells.Item(126,3).Value = 0
$ws.Cells.Item(126,4).Value = 1
$ws.Cells.Item(126,5).Value = 36
$ws.Cells.Item(126,6).Value = 0
$ws.Cells.Item(126,7).Value = 0
$ws.Cells.Item(126,8).Value = 2
$ws.Cells.Item(127,1).Value = "Macao"
$ws.Cells.Item(127,2).Value = 34
$ws.Cells.Item(127,3).Value = 1
$ws.Cells.Item(127,4).Value = 10
$ws.Cells.Item(127,5).Value = 24
$ws.Cells.Item(127,6).Value = 0
$ws.Cells.Item(127,7).Value = 0
$ws.Cells.Item(127,8).Value = 0
$ws.Cells.Item(128,1).Value = "Monaco"
$ws.Cells.Item(128,2).Value = 33
$ws.Cells.Item(128,3).Value = 0
$ws.Cells.Item(128,4).Value = 1
$ws.Cells.Item(128,5).Value = 32
$ws.Cells.Item(128,6).Value = 0
$ws.Cells.Item(128,7).Value = 0
$ws.Cells.Item(128,8).Value = 0
$ws.Cells.Item(129,1).Value = "Guam"
$ws.Cells.Item(129,2).Value = 32
$ws.Cells.Item(129,3).Value = 0
$ws.Cells.Item(129,4).Value = 0
$ws.Cells.Item(129,5).Value = 31
$ws.Cells.Item(129,6).Value = 0
$ws.Cells.Item(129,7).Value = 0
$ws.Cells.Item(129,8).Value = 1
$ws.Cells.Item(130,1).Value = "Kenia"
$ws.Cells.Item(130,2).Value = 31
$ws.Cells.Item(130,3).Value = 0
$ws.Cells.Item(130,4).Value = 1
$ws.Cells.Item(130,5).Value = 29
$ws.Cells.Item(130,6).Value = 0
$ws.Cells.Item(130,7).Value = 0
$ws.Cells.Item(130,8).Value = 1
$ws.Cells.Item(131,1).Value = "Polinesia Francesa"
$ws.Cells.Item(131,2).Value = 30
$ws.Cells.Item(131,3).Value = 0
$ws.Cells.Item(131,4).Value = 0
$ws.Cells.Item(131,5).Value = 30
$ws.Cells.Item(131,6).Value = 0
$ws.Cells.Item(131,7).Value = 0
$ws.Cells.Item(131,8).Value = 0
$ws.Cells.Item(132,1).Value = "Isla de Man"
$ws.Cells.Item(132,2).Value = 29
$ws.Cells.Item(132,3).Value = 3
$ws.Cells.Item(132,4).Value = 0
$ws.Cells.Item(132,5).Value = 29
$ws.Cells.Item(132,6).Value = 0
$ws.Cells.Item(132,7).Value = 0
$ws.Cells.Item(132,8).Value = 0
$ws.Cells.Item(133,1).Value = "Aruba"
$ws.Cells.Item(133,2).Value = 28
$ws.Cells.Item(133,3).Value = 0
$ws.Cells.Item(133,4).Value = 1
$ws.Cells.Item(133,5).Value = 27
$ws.Cells.Item(133,6).Value = 0
$ws.Cells.Item(133,7).Value = 0
$ws.Cells.Item(133,8).Value = 0
$ws.Cells.Item(134,1).Value = "Guayana Francesa"
$ws.Cells.Item(134,2).Value = 28
$ws.Cells.Item(134,3).Value = 0
$ws.Cells.Item(134,4).Value = 6
$ws.Cells.Item(134,5).Value = 22
$ws.Cells.Item(134,6).Value = 0
$ws.Cells.Item(134,7).Value = 0
$ws.Cells.Item(134,8).Value = 0
$ws.Cells.Item(135,1).Value = "Jamaica"
$ws.Cells.Item(135,2).Value = 26
$ws.Cells.Item(135,3).Value = 0
$ws.Cells.Item(135,4).Value = 2
$ws.Cells.Item(135,5).Value = 23
$ws.Cells.Item(135,6).Value = 0
$ws.Cells.Item(135,7).Value = 0
$ws.Cells.Item(135,8).Value = 1
$ws.Cells.Item(136,1).Value = "Togo"
$ws.Cells.Item(136,2).Value = 25
$ws.Cells.Item(136,3).Value = 1
$ws.Cells.Item(136,4).Value = 1
$ws.Cells.Item(136,5).Value = 24
$ws.Cells.Item(136,6).Value = 0
$ws.Cells.Item(136,7).Value = 0
$ws.Cells.Item(136,8).Value = 0
$ws.Cells.Item(137,1).Value = "Guatemala"
$ws.Cells.Item(137,2).Value = 25
$ws.Cells.Item(137,3).Value = 0
$ws.Cells.Item(137,4).Value = 4
$ws.Cells.Item(137,5).Value = 20
$ws.Cells.Item(137,6).Value = 0
$ws.Cells.Item(137,7).Value = 0
$ws.Cells.Item(137,8).Value = 1
$ws.Cells.Item(138,1).Value = "Barbados"
$ws.Cells.Item(138,2).Value = 24
$ws.Cells.Item(138,3).Value = 0
$ws.Cells.Item(138,4).Value = 0
$ws.Cells.Item(138,5).Value = 24
$ws.Cells.Item(138,6).Value = 0
$ws.Cells.Item(138,7).Value = 0
$ws.Cells.Item(138,8).Value = 0
$ws.Cells.Item(139,1).Value = "Madagascar"
$ws.Cells.Item(139,2).Value = 24
$ws.Cells.Item(139,3).Value = 1
$ws.Cells.Item(139,4).Value = 0
$ws.Cells.Item(139,5).Value = 24
$ws.Cells.Item(139,6).Value = 0
$ws.Cells.Item(139,7).Value = 0
$ws.Cells.Item(139,8).Value = 0
$ws.Cells.Item(140,1).Value = "Zambia"
$ws.Cells.Item(140,2).Value = 22
$ws.Cells.Item(140,3).Value = 6
$ws.Cells.Item(140,4).Value = 0
$ws.Cells.Item(140,5).Value = 22
$ws.Cells.Item(140,6).Value = 0
$ws.Cells.Item(140,7).Value = 0
$ws.Cells.Item(140,8).Value = 0
$ws.Cells.Item(141,1).Value = "Uganda"
$ws.Cells.Item(141,2).Value = 18
$ws.Cells.Item(141,3).Value = 4
$ws.Cells.Item(141,4).Value = 0
$ws.Cells.Item(141,5).Value = 18
$ws.Cells.Item(141,6).Value = 0
$ws.Cells.Item(141,7).Value = 0
$ws.Cells.Item(141,8).Value = 0
$ws.Cells.Item(142,1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(142,2).Value = 17
$ws.Cells.Item(142,3).Value = 0
$ws.Cells.Item(142,4).Value = 0
$ws.Cells.Item(142,5).Value = 17
$ws.Cells.Item(142,6).Value = 0
$ws.Cells.Item(142,7).Value = 0
$ws.Cells.Item(142,8).Value = 0
$ws.Cells.Item(143,1).Value = "Etiopia"
$ws.Cells.Item(143,2).Value = 16
$ws.Cells.Item(143,3).Value = 4
$ws.Cells.Item(143,4).Value = 0
$ws.Cells.Item(143,5).Value = 16
$ws.Cells.Item(143,6).Value = 0
$ws.Cells.Item(143,7).Value = 0
$ws.Cells.Item(143,8).Value = 0
$ws.Cells.Item(144,1).Value = "Nueva Caledonia"
$ws.Cells.Item(144,2).Value = 15
$ws.Cells.Item(144,3).Value = 1
$ws.Cells.Item(144,4).Value = 0
$ws.Cells.Item(144,5).Value = 15
$ws.Cells.Item(144,6).Value = 0
$ws.Cells.Item(144,7).Value = 0
$ws.Cells.Item(144,8).Value = 0
$ws.Cells.Item(145,1).Value = "Bermudas"
$ws.Cells.Item(145,2).Value = 15
$ws.Cells.Item(145,3).Value = 0
$ws.Cells.Item(145,4).Value = 2
$ws.Cells.Item(145,5).Value = 13
$ws.Cells.Item(145,6).Value = 0
$ws.Cells.Item(145,7).Value = 0
$ws.Cells.Item(145,8).Value = 0
$ws.Cells.Item(146,1).Value = "Maldivas"
$ws.Cells.Item(146,2).Value = 14
$ws.Cells.Item(146,3).Value = 1
$ws.Cells.Item(146,4).Value = 9
$ws.Cells.Item(146,5).Value = 5
$ws.Cells.Item(146,6).Value = 0
$ws.Cells.Item(146,7).Value = 0
$ws.Cells.Item(146,8).Value = 0
$ws.Cells.Item(147,1).Value = "El Salvador"
$ws.Cells.Item(147,2).Value = 13
$ws.Cells.Item(147,3).Value = 0
$ws.Cells.Item(147,4).Value = 0
$ws.Cells.Item(147,5).Value = 13
$ws.Cells.Item(147,6).Value = 0
$ws.Cells.Item(147,7).Value = 0
$ws.Cells.Item(147,8).Value = 0
$ws.Cells.Item(148,1).Value = "Tanzania"
$ws.Cells.Item(148,2).Value = 13
$ws.Cells.Item(148,3).Value = 0
$ws.Cells.Item(148,4).Value = 1
$ws.Cells.Item(148,5).Value = 12
$ws.Cells.Item(148,6).Value = 0
$ws.Cells.Item(148,7).Value = 0
$ws.Cells.Item(148,8).Value = 0
$ws.Cells.Item(149,1).Value = "Republica de Yibuti"
$ws.Cells.Item(149,2).Value = 12
$ws.Cells.Item(149,3).Value = 1
$ws.Cells.Item(149,4).Value = 0
$ws.Cells.Item(149,5).Value = 12
$ws.Cells.Item(149,6).Value = 0
$ws.Cells.Item(149,7).Value = 0
$ws.Cells.Item(149,8).Value = 0
$ws.Cells.Item(150,1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(150,2).Value = 12
$ws.Cells.Item(150,3).Value = 0
$ws.Cells.Item(150,4).Value = 0
$ws.Cells.Item(150,5).Value = 12
$ws.Cells.Item(150,6).Value = 0
$ws.Cells.Item(150,7).Value = 0
$ws.Cells.Item(150,8).Value = 0
$ws.Cells.Item(151,1).Value = "Mongolia"
$ws.Cells.Item(151,2).Value = 11
$ws.Cells.Item(151,3).Value = 0
$ws.Cells.Item(151,4).Value = 0
$ws.Cells.Item(151,5).Value = 11
$ws.Cells.Item(151,6).Value = 0
$ws.Cells.Item(151,7).Value = 0
$ws.Cells.Item(151,8).Value = 0
$ws.Cells.Item(152,1).Value = "Dominica"
$ws.Cells.Item(152,2).Value = 11
$ws.Cells.Item(152,3).Value = 0
$ws.Cells.Item(152,4).Value = 0
$ws.Cells.Item(152,5).Value = 11
$ws.Cells.Item(152,6).Value = 0
$ws.Cells.Item(152,7).Value = 0
$ws.Cells.Item(152,8).Value = 0
$ws.Cells.Item(153,1).Value = "Mali"
$ws.Cells.Item(153,2).Value = 11
$ws.Cells.Item(153,3).Value = 7
$ws.Cells.Item(153,4).Value = 0
$ws.Cells.Item(153,5).Value = 11
$ws.Cells.Item(153,6).Value = 0
$ws.Cells.Item(153,7).Value = 0
$ws.Cells.Item(153,8).Value = 0
$ws.Cells.Item(154,1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(154,2).Value = 11
$ws.Cells.Item(154,3).Value = 0
$ws.Cells.Item(154,4).Value = 0
$ws.Cells.Item(154,5).Value = 11
$ws.Cells.Item(154,6).Value = 0
$ws.Cells.Item(154,7).Value = 0
$ws.Cells.Item(154,8).Value = 0
$ws.Cells.Item(155,1).Value = "Niger"
$ws.Cells.Item(155,2).Value = 10
$ws.Cells.Item(155,3).Value = 0
$ws.Cells.Item(155,4).Value = 0
$ws.Cells.Item(155,5).Value = 9
$ws.Cells.Item(155,6).Value = 0
$ws.Cells.Item(155,7).Value = 0
$ws.Cells.Item(155,8).Value = 1
$ws.Cells.Item(156,1).Value = "Groenlandia"
$ws.Cells.Item(156,2).Value = 10
$ws.Cells.Item(156,3).Value = 4
$ws.Cells.Item(156,4).Value = 2
$ws.Cells.Item(156,5).Value = 8
$ws.Cells.Item(156,6).Value = 0
$ws.Cells.Item(156,7).Value = 0
$ws.Cells.Item(156,8).Value = 0
$ws.Cells.Item(157,1).Value = "Suazilandia"
$ws.Cells.Item(157,2).Value = 9
$ws.Cells.Item(157,3).Value = 3
$ws.Cells.Item(157,4).Value = 0
$ws.Cells.Item(157,5).Value = 9
$ws.Cells.Item(157,6).Value = 0
$ws.Cells.Item(157,7).Value = 0
$ws.Cells.Item(157,8).Value = 0
$ws.Cells.Item(158,1).Value = "Bahamas"
$ws.Cells.Item(158,2).Value = 9
$ws.Cells.Item(158,3).Value = 0
$ws.Cells.Item(158,4).Value = 1
$ws.Cells.Item(158,5).Value = 8
$ws.Cells.Item(158,6).Value = 0
$ws.Cells.Item(158,7).Value = 0
$ws.Cells.Item(158,8).Value = 0
$ws.Cells.Item(159,1).Value = "Surinam"
$ws.Cells.Item(159,2).Value = 8
$ws.Cells.Item(159,3).Value = 0
$ws.Cells.Item(159,4).Value = 0
$ws.Cells.Item(159,5).Value = 8
$ws.Cells.Item(159,6).Value = 0
$ws.Cells.Item(159,7).Value = 0
$ws.Cells.Item(159,8).Value = 0
$ws.Cells.Item(160,1).Value = "Guinea"
$ws.Cells.Item(160,2).Value = 8
$ws.Cells.Item(160,3).Value = 4
$ws.Cells.Item(160,4).Value = 0
$ws.Cells.Item(160,5).Value = 8
$ws.Cells.Item(160,6).Value = 0
$ws.Cells.Item(160,7).Value = 0
$ws.Cells.Item(160,8).Value = 0
$ws.Cells.Item(161,1).Value = "Haiti"
$ws.Cells.Item(161,2).Value = 8
$ws.Cells.Item(161,3).Value = 0
$ws.Cells.Item(161,4).Value = 0
$ws.Cells.Item(161,5).Value = 8
$ws.Cells.Item(161,6).Value = 0
$ws.Cells.Item(161,7).Value = 0
$ws.Cells.Item(161,8).Value = 0
$ws.Cells.Item(162,1).Value = "Islas Caimanes"
$ws.Cells.Item(162,2).Value = 8
$ws.Cells.Item(162,3).Value = 0
$ws.Cells.Item(162,4).Value = 0
$ws.Cells.Item(162,5).Value = 7
$ws.Cells.Item(162,6).Value = 0
$ws.Cells.Item(162,7).Value = 0
$ws.Cells.Item(162,8).Value = 1
$ws.Cells.Item(163,1).Value = "Namibia"
$ws.Cells.Item(163,2).Value = 8
$ws.Cells.Item(163,3).Value = 0
$ws.Cells.Item(163,4).Value = 2
$ws.Cells.Item(163,5).Value = 6
$ws.Cells.Item(163,6).Value = 0
$ws.Cells.Item(163,7).Value = 0
$ws.Cells.Item(163,8).Value = 0
$ws.Cells.Item(164,1).Value = "Seychelles"
$ws.Cells.Item(164,2).Value = 7
$ws.Cells.Item(164,3).Value = 0
$ws.Cells.Item(164,4).Value = 0
$ws.Cells.Item(164,5).Value = 7
$ws.Cells.Item(164,6).Value = 0
$ws.Cells.Item(164,7).Value = 0
$ws.Cells.Item(164,8).Value = 0
$ws.Cells.Item(165,1).Value = "Antigua y Barbuda"
$ws.Cells.Item(165,2).Value = 7
$ws.Cells.Item(165,3).Value = 0
$ws.Cells.Item(165,4).Value = 0
$ws.Cells.Item(165,5).Value = 7
$ws.Cells.Item(165,6).Value = 0
$ws.Cells.Item(165,7).Value = 0
$ws.Cells.Item(165,8).Value = 0
$ws.Cells.Item(166,1).Value = "Mozambique"
$ws.Cells.Item(166,2).Value = 7
$ws.Cells.Item(166,3).Value = 0
$ws.Cells.Item(166,4).Value = 0
$ws.Cells.Item(166,5).Value = 7
$ws.Cells.Item(166,6).Value = 0
$ws.Cells.Item(166,7).Value = 0
$ws.Cells.Item(166,8).Value = 0
$ws.Cells.Item(167,1).Value = "Granada"
$ws.Cells.Item(167,2).Value = 7
$ws.Cells.Item(167,3).Value = 0
$ws.Cells.Item(167,4).Value = 0
$ws.Cells.Item(167,5).Value = 7
$ws.Cells.Item(167,6).Value = 0
$ws.Cells.Item(167,7).Value = 0
$ws.Cells.Item(167,8).Value = 0
$ws.Cells.Item(168,1).Value = "Gabon"
$ws.Cells.Item(168,2).Value = 7
$ws.Cells.Item(168,3).Value = 0
$ws.Cells.Item(168,4).Value = 0
$ws.Cells.Item(168,5).Value = 6
$ws.Cells.Item(168,6).Value = 0
$ws.Cells.Item(168,7).Value = 0
$ws.Cells.Item(168,8).Value = 1
$ws.Cells.Item(169,1).Value = "Curazao"
$ws.Cells.Item(169,2).Value = 7
$ws.Cells.Item(169,3).Value = 0
$ws.Cells.Item(169,4).Value = 2
$ws.Cells.Item(169,5).Value = 4
$ws.Cells.Item(169,6).Value = 0
$ws.Cells.Item(169,7).Value = 0
$ws.Cells.Item(169,8).Value = 1
$ws.Cells.Item(170,1).Value = "Laos"
$ws.Cells.Item(170,2).Value = 6
$ws.Cells.Item(170,3).Value = 0
$ws.Cells.Item(170,4).Value = 0
$ws.Cells.Item(170,5).Value = 6
$ws.Cells.Item(170,6).Value = 0
$ws.Cells.Item(170,7).Value = 0
$ws.Cells.Item(170,8).Value = 0
$ws.Cells.Item(171,1).Value = "Benin"
$ws.Cells.Item(171,2).Value = 6
$ws.Cells.Item(171,3).Value = 0
$ws.Cells.Item(171,4).Value = 0
$ws.Cells.Item(171,5).Value = 6
$ws.Cells.Item(171,6).Value = 0
$ws.Cells.Item(171,7).Value = 0
$ws.Cells.Item(171,8).Value = 0
$ws.Cells.Item(172,1).Value = "Eritrea"
$ws.Cells.Item(172,2).Value = 6
$ws.Cells.Item(172,3).Value = 0
$ws.Cells.Item(172,4).Value = 0
$ws.Cells.Item(172,5).Value = 6
$ws.Cells.Item(172,6).Value = 0
$ws.Cells.Item(172,7).Value = 0
$ws.Cells.Item(172,8).Value = 0
$ws.Cells.Item(173,1).Value = "Siria"
$ws.Cells.Item(173,2).Value = 5
$ws.Cells.Item(173,3).Value = 0
$ws.Cells.Item(173,4).Value = 0
$ws.Cells.Item(173,5).Value = 5
$ws.Cells.Item(173,6).Value = 0
$ws.Cells.Item(173,7).Value = 0
$ws.Cells.Item(173,8).Value = 0
$ws.Cells.Item(174,1).Value = "Montserrat"
$ws.Cells.Item(174,2).Value = 5
$ws.Cells.Item(174,3).Value = 0
$ws.Cells.Item(174,4).Value = 0
$ws.Cells.Item(174,5).Value = 5
$ws.Cells.Item(174,6).Value = 0
$ws.Cells.Item(174,7).Value = 0
$ws.Cells.Item(174,8).Value = 0
$ws.Cells.Item(175,1).Value = "Fiyi"
$ws.Cells.Item(175,2).Value = 5
$ws.Cells.Item(175,3).Value = 0
$ws.Cells.Item(175,4).Value = 0
$ws.Cells.Item(175,5).Value = 5
$ws.Cells.Item(175,6).Value = 0
$ws.Cells.Item(175,7).Value = 0
$ws.Cells.Item(175,8).Value = 0
$ws.Cells.Item(176,1).Value = "Birmania"
$ws.Cells.Item(176,2).Value = 5
$ws.Cells.Item(176,3).Value = 0
$ws.Cells.Item(176,4).Value = 0
$ws.Cells.Item(176,5).Value = 5
$ws.Cells.Item(176,6).Value = 0
$ws.Cells.Item(176,7).Value = 0
$ws.Cells.Item(176,8).Value = 0
$ws.Cells.Item(177,1).Value = "Cabo Verde"
$ws.Cells.Item(177,2).Value = 5
$ws.Cells.Item(177,3).Value = 0
$ws.Cells.Item(177,4).Value = 0
$ws.Cells.Item(177,5).Value = 4
$ws.Cells.Item(177,6).Value = 0
$ws.Cells.Item(177,7).Value = 0
$ws.Cells.Item(177,8).Value = 1
$ws.Cells.Item(178,1).Value = "Zimbabue"
$ws.Cells.Item(178,2).Value = 5
$ws.Cells.Item(178,3).Value = 2
$ws.Cells.Item(178,4).Value = 0
$ws.Cells.Item(178,5).Value = 4
$ws.Cells.Item(178,6).Value = 0
$ws.Cells.Item(178,7).Value = 0
$ws.Cells.Item(178,8).Value = 1
$ws.Cells.Item(179,1).Value = "Guyana"
$ws.Cells.Item(179,2).Value = 5
$ws.Cells.Item(179,3).Value = 0
$ws.Cells.Item(179,4).Value = 0
$ws.Cells.Item(179,5).Value = 4
$ws.Cells.Item(179,6).Value = 0
$ws.Cells.Item(179,7).Value = 0
$ws.Cells.Item(179,8).Value = 1
$ws.Cells.Item(180,1).Value = "Congo"
$ws.Cells.Item(180,2).Value = 4
$ws.Cells.Item(180,3).Value = 0
$ws.Cells.Item(180,4).Value = 0
$ws.Cells.Item(180,5).Value = 4
$ws.Cells.Item(180,6).Value = 0
$ws.Cells.Item(180,7).Value = 0
$ws.Cells.Item(180,8).Value = 0
$ws.Cells.Item(181,1).Value = "Santa Sede"
$ws.Cells.Item(181,2).Value = 4
$ws.Cells.Item(181,3).Value = 0
$ws.Cells.Item(181,4).Value = 0
$ws.Cells.Item(181,5).Value = 4
$ws.Cells.Item(181,6).Value = 0
$ws.Cells.Item(181,7).Value = 0
$ws.Cells.Item(181,8).Value = 0
$ws.Cells.Item(182,1).Value = "Angola"
$ws.Cells.Item(182,2).Value = 4
$ws.Cells.Item(182,3).Value = 0
$ws.Cells.Item(182,4).Value = 0
$ws.Cells.Item(182,5).Value = 4
$ws.Cells.Item(182,6).Value = 0
$ws.Cells.Item(182,7).Value = 0
$ws.Cells.Item(182,8).Value = 0
$ws.Cells.Item(183,1).Value = "Nepal"
$ws.Cells.Item(183,2).Value = 4
$ws.Cells.Item(183,3).Value = 1
$ws.Cells.Item(183,4).Value = 1
$ws.Cells.Item(183,5).Value = 3
$ws.Cells.Item(183,6).Value = 0
$ws.Cells.Item(183,7).Value = 0
$ws.Cells.Item(183,8).Value = 0
$ws.Cells.Item(184,1).Value = "Republica del Chad"
$ws.Cells.Item(184,2).Value = 3
$ws.Cells.Item(184,3).Value = 0
$ws.Cells.Item(184,4).Value = 0
$ws.Cells.Item(184,5).Value = 3
$ws.Cells.Item(184,6).Value = 0
$ws.Cells.Item(184,7).Value = 0
$ws.Cells.Item(184,8).Value = 0
$ws.Cells.Item(185,1).Value = "Republica de Africa Central"
$ws.Cells.Item(185,2).Value = 3
$ws.Cells.Item(185,3).Value = 0
$ws.Cells.Item(185,4).Value = 0
$ws.Cells.Item(185,5).Value = 3
$ws.Cells.Item(185,6).Value = 0
$ws.Cells.Item(185,7).Value = 0
$ws.Cells.Item(185,8).Value = 0
$ws.Cells.Item(186,1).Value = "Somalia"
$ws.Cells.Item(186,2).Value = 3
$ws.Cells.Item(186,3).Value = 1
$ws.Cells.Item(186,4).Value = 0
$ws.Cells.Item(186,5).Value = 3
$ws.Cells.Item(186,6).Value = 0
$ws.Cells.Item(186,7).Value = 0
$ws.Cells.Item(186,8).Value = 0
$ws.Cells.Item(187,1).Value = "Butan"
$ws.Cells.Item(187,2).Value = 3
$ws.Cells.Item(187,3).Value = 1
$ws.Cells.Item(187,4).Value = 0
$ws.Cells.Item(187,5).Value = 3
$ws.Cells.Item(187,6).Value = 0
$ws.Cells.Item(187,7).Value = 0
$ws.Cells.Item(187,8).Value = 0
$ws.Cells.Item(188,1).Value = "Mauritania"
$ws.Cells.Item(188,2).Value = 3
$ws.Cells.Item(188,3).Value = 0
$ws.Cells.Item(188,4).Value = 0
$ws.Cells.Item(188,5).Value = 3
$ws.Cells.Item(188,6).Value = 0
$ws.Cells.Item(188,7).Value = 0
$ws.Cells.Item(188,8).Value = 0
$ws.Cells.Item(189,1).Value = "Liberia"
$ws.Cells.Item(189,2).Value = 3
$ws.Cells.Item(189,3).Value = 0
$ws.Cells.Item(189,4).Value = 0
$ws.Cells.Item(189,5).Value = 3
$ws.Cells.Item(189,6).Value = 0
$ws.Cells.Item(189,7).Value = 0
$ws.Cells.Item(189,8).Value = 0
$ws.Cells.Item(190,1).Value = "San Bartolome"
$ws.Cells.Item(190,2).Value = 3
$ws.Cells.Item(190,3).Value = 0
$ws.Cells.Item(190,4).Value = 0
$ws.Cells.Item(190,5).Value = 3
$ws.Cells.Item(190,6).Value = 0
$ws.Cells.Item(190,7).Value = 0
$ws.Cells.Item(190,8).Value = 0
$ws.Cells.Item(191,1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(191,2).Value = 3
$ws.Cells.Item(191,3).Value = 0
$ws.Cells.Item(191,4).Value = 0
$ws.Cells.Item(191,5).Value = 3
$ws.Cells.Item(191,6).Value = 0
$ws.Cells.Item(191,7).Value = 0
$ws.Cells.Item(191,8).Value = 0
$ws.Cells.Item(192,1).Value = "Gambia"
$ws.Cells.Item(192,2).Value = 3
$ws.Cells.Item(192,3).Value = 0
$ws.Cells.Item(192,4).Value = 0
$ws.Cells.Item(192,5).Value = 2
$ws.Cells.Item(192,6).Value = 0
$ws.Cells.Item(192,7).Value = 0
$ws.Cells.Item(192,8).Value = 1
$ws.Cells.Item(193,1).Value = "Sudan"
$ws.Cells.Item(193,2).Value = 3
$ws.Cells.Item(193,3).Value = 0
$ws.Cells.Item(193,4).Value = 0
$ws.Cells.Item(193,5).Value = 2
$ws.Cells.Item(193,6).Value = 0
$ws.Cells.Item(193,7).Value = 0
$ws.Cells.Item(193,8).Value = 1
$ws.Cells.Item(194,1).Value = "Santa Lucia"
$ws.Cells.Item(194,2).Value = 3
$ws.Cells.Item(194,3).Value = 0
$ws.Cells.Item(194,4).Value = 1
$ws.Cells.Item(194,5).Value = 2
$ws.Cells.Item(194,6).Value = 0
$ws.Cells.Item(194,7).Value = 0
$ws.Cells.Item(194,8).Value = 0
$ws.Cells.Item(195,1).Value = "Belice"
$ws.Cells.Item(195,2).Value = 2
$ws.Cells.Item(195,3).Value = 0
$ws.Cells.Item(195,4).Value = 0
$ws.Cells.Item(195,5).Value = 2
$ws.Cells.Item(195,6).Value = 0
$ws.Cells.Item(195,7).Value = 0
$ws.Cells.Item(195,8).Value = 0
$ws.Cells.Item(196,1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(196,2).Value = 2
$ws.Cells.Item(196,3).Value = 0
$ws.Cells.Item(196,4).Value = 0
$ws.Cells.Item(196,5).Value = 2
$ws.Cells.Item(196,6).Value = 0
$ws.Cells.Item(196,7).Value = 0
$ws.Cells.Item(196,8).Value = 0
$ws.Cells.Item(197,1).Value = "Anguila"
$ws.Cells.Item(197,2).Value = 2
$ws.Cells.Item(197,3).Value = 0
$ws.Cells.Item(197,4).Value = 0
$ws.Cells.Item(197,5).Value = 2
$ws.Cells.Item(197,6).Value = 0
$ws.Cells.Item(197,7).Value = 0
$ws.Cells.Item(197,8).Value = 0
$ws.Cells.Item(198,1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(198,2).Value = 2
$ws.Cells.Item(198,3).Value = 0
$ws.Cells.Item(198,4).Value = 0
$ws.Cells.Item(198,5).Value = 2
$ws.Cells.Item(198,6).Value = 0
$ws.Cells.Item(198,7).Value = 0
$ws.Cells.Item(198,8).Value = 0
$ws.Cells.Item(199,1).Value = "Guinea-Bisau"
$ws.Cells.Item(199,2).Value = 2
$ws.Cells.Item(199,3).Value = 0
$ws.Cells.Item(199,4).Value = 0
$ws.Cells.Item(199,5).Value = 2
$ws.Cells.Item(199,6).Value = 0
$ws.Cells.Item(199,7).Value = 0
$ws.Cells.Item(199,8).Value = 0
$ws.Cells.Item(200,1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(200,2).Value = 2
$ws.Cells.Item(200,3).Value = 0
$ws.Cells.Item(200,4).Value = 0
$ws.Cells.Item(200,5).Value = 2
$ws.Cells.Item(200,6).Value = 0
$ws.Cells.Item(200,7).Value = 0
$ws.Cells.Item(200,8).Value = 0
$ws.Cells.Item(201,1).Value = "Nicaragua"
$ws.Cells.Item(201,2).Value = 2
$ws.Cells.Item(201,3).Value = 0
$ws.Cells.Item(201,4).Value = 0
$ws.Cells.Item(201,5).Value = 1
$ws.Cells.Item(201,6).Value = 0
$ws.Cells.Item(201,7).Value = 1
$ws.Cells.Item(201,8).Value = 1
$ws.Cells.Item(202,1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(202,2).Value = 1
$ws.Cells.Item(202,3).Value = 0
$ws.Cells.Item(202,4).Value = 0
$ws.Cells.Item(202,5).Value = 1
$ws.Cells.Item(202,6).Value = 0
$ws.Cells.Item(202,7).Value = 0
$ws.Cells.Item(202,8).Value = 0
$ws.Cells.Item(203,1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(203,2).Value = 1
$ws.Cells.Item(203,3).Value = 0
$ws.Cells.Item(203,4).Value = 0
$ws.Cells.Item(203,5).Value = 1
$ws.Cells.Item(203,6).Value = 0
$ws.Cells.Item(203,7).Value = 0
$ws.Cells.Item(203,8).Value = 0
$ws.Cells.Item(204,1).Value = "Timor Oriental"
$ws.Cells.Item(204,2).Value = 1
$ws.Cells.Item(204,3).Value = 0
$ws.Cells.Item(204,4).Value = 0
$ws.Cells.Item(204,5).Value = 1
$ws.Cells.Item(204,6).Value = 0
$ws.Cells.Item(204,7).Value = 0
$ws.Cells.Item(204,8).Value = 0
$ws.Cells.Item(205,1).Value = "Libia"
$ws.Cells.Item(205,2).Value = 1
$ws.Cells.Item(205,3).Value = 0
$ws.Cells.Item(205,4).Value = 0
$ws.Cells.Item(205,5).Value = 1
$ws.Cells.Item(205,6).Value = 0
$ws.Cells.Item(205,7).Value = 0
$ws.Cells.Item(205,8).Value = 0

Write-Host "Applied country data update"
